# Replaces the placeholder zero rows (6-19) with real recorded event data and
# refreshes the already-recorded TR rows (2-5) with the corrected
# RelativeTime/AbsoluteTime/Difference values (new participant event export).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row layout: A=TR, B=RelativeTime, C=AbsoluteTime, D=Difference, E=Description
$rows = @(
    @{ Row = 2;  A = 1;  B = 1.9997048999648541;  C = 2018076.0184128999; D = 2018078.0202494999; E = "Active Stimuli" },
    @{ Row = 3;  A = 2;  B = 3.9987212000414729;  C = 2018076.0185197999; D = 2018080.0192658;     E = "Active Stimuli" },
    @{ Row = 4;  A = 3;  B = 5.998412000015378;   C = 2018076.0186335;    D = 2018082.0189566;     E = "Active Stimuli" },
    @{ Row = 5;  A = 4;  B = 7.9983878000639379;  C = 2018076.0186334001; D = 2018084.0189324;     E = "Active Stimuli" },
    @{ Row = 6;  A = 5;  B = 9.9993088999763131;  C = 2018076.0186073;    D = 2018086.0198535;     E = "Active Stimuli" },
    @{ Row = 7;  A = 6;  B = 11.998830100055784;  C = 2018076.0186324001; D = 2018088.0193747;     E = "Active Stimuli" },
    @{ Row = 8;  A = 7;  B = 13.999511400004849;  C = 2018076.0186307;    D = 2018090.020056;      E = "Active Stimuli" },
    @{ Row = 9;  A = 8;  B = 15.998318700119853;  C = 2018076.0186332001; D = 2018092.0188633001;  E = "Active Stimuli" },
    @{ Row = 10; A = 9;  B = 17.998543200083077;  C = 2018076.0186333;    D = 2018094.0190878001;  E = "Active Stimuli" },
    @{ Row = 11; A = 10; B = 19.998396700015292;  C = 2018076.0186331;    D = 2018096.0189413;     E = "Active Stimuli" },
    @{ Row = 12; A = 11; B = 21.998372799949721;  C = 2018076.0186329;    D = 2018098.0189173999;  E = "Active Stimuli" },
    @{ Row = 13; A = 12; B = 23.998295200057328;  C = 2018076.0186334001; D = 2018100.0188398;     E = "Active Stimuli" },
    @{ Row = 14; A = 13; B = 25.998398100025952;  C = 2018076.0186331;    D = 2018102.0189427;     E = "Active Stimuli" },
    @{ Row = 15; A = 14; B = 27.998278399929404;  C = 2018076.0186333;    D = 2018104.0188229999;  E = "Active Stimuli" },
    @{ Row = 16; A = 15; B = 29.998560100095347;  C = 2018076.0186313;    D = 2018106.0191047001;  E = "Active Stimuli" },
    @{ Row = 17; A = 16; B = 31.998352100141346;  C = 2018076.0186328001; D = 2018108.0188967001;  E = "Active Stimuli" },
    @{ Row = 18; A = 17; B = 33.998800799949095;  C = 2018076.018499;     D = 2018110.0193453999;  E = "Cross" },
    @{ Row = 19; A = 18; B = 35.998359200078994;  C = 2018076.0186292001; D = 2018112.0189038001;  E = "Cross" }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.A
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D
    $ws.Cells.Item($r.Row, 5).Value = $r.E
}
